# "updating the clash file from bat" — refresh the Raids sheet with the
# latest clan-member snapshot: trophy counts bumped for several members,
# and three pairs of rows swap places (the underlying roster re-sorted by
# the bot, but row styling/rank stayed put, so tag/trophy/name/DR values
# had to be re-seated row by row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Raids")

# Trophy (column B) refresh - no other data moved for these rows.
$ws.Range("B12").Value = 2922
$ws.Range("B15").Value = 1412
$ws.Range("B19").Value = 1183
$ws.Range("B29").Value = 864
$ws.Range("B31").Value = 1531
$ws.Range("B33").Value = 1192

# Row 34 now holds the member previously shown on row 39 (SUPoT).
$ws.Range("A34").Value = "#QQL28Y2UL"
$ws.Range("B34").Value = 1112
$ws.Range("D34").Value = "SUPoT"
$ws.Range("H34").Value = 56

# Rows 36 and 37 swap (Zodiac <-> Apollo).
$ws.Range("A36").Value = "#QV8RY9UC8"
$ws.Range("B36").Value = 817
$ws.Range("D36").Value = "Apollo"

$ws.Range("A37").Value = "#LPCLQUCCY"
$ws.Range("B37").Value = 1013
$ws.Range("D37").Value = "Zodiac"

# Row 39 now holds the member previously shown on row 34 (Kingsman).
$ws.Range("A39").Value = "#QUGYGY88C"
$ws.Range("B39").Value = 941
$ws.Range("D39").Value = "Kingsman"
$ws.Range("H39").Value = 0

# Rows 43 and 44 swap (Death1wolf <-> DaddyChill).
$ws.Range("A43").Value = "#G0LJCVR2P"
$ws.Range("B43").Value = 1051
$ws.Range("D43").Value = "DaddyChill"

$ws.Range("A44").Value = "#LGCVY0L9P"
$ws.Range("B44").Value = 935
$ws.Range("D44").Value = "Death1wolf"

# Row 46 now holds the member previously shown on row 43/44 chain (Kukoshibo).
$ws.Range("A46").Value = "#QLUV29GGJ"
$ws.Range("B46").Value = 994
$ws.Range("D46").Value = "Kukoshibo"

# Name column (D) got a bit wider on re-save once Excel recomputed the
# best-fit width for the refreshed names.
$ws.Columns.Item(4).AutoFit() | Out-Null
